$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 14: "Add Secret Sound for Cave Secret + Compose" - Effort (D14) increased from 2 to 3 (task completed)
$ws.Range("D14").Value = 3

# Row 20: new task "Implementation of First Boss ("Witch Tree") with Squirel"
$ws.Range("A20").Value = 'Implementation of First Boss ("Witch Tree") with Squirel'
$ws.Range("B20").Value = 10
$ws.Range("C20").Value = 10
$ws.Range("F20").Value = "Sascha"

# Row 21: new task "Combatsystem for Axe"
$ws.Range("A21").Value = "Combatsystem for Axe"

# Update the active selection to I20 to reflect where the user left off
$ws.Range("I20").Select() | Out-Null
